$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '42.982.23'
$ws.Range("E2").Value = '  +4.34%  '
$ws.Range("D3").Value = '2.239.15'
$ws.Range("E3").Value = '  +2.92%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.619'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.80%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '76.05'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +8.02%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.617'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.93'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0933'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.63'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.98%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.98'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.61%  '
$ws.Range("E14").Value = '  +0.48%  '
$ws.Range("D15").Value = '2.555.35'
$ws.Range("E15").Value = '  +2.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.69'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.31%  '
$ws.Range("D17").Value = '2.237.41'
$ws.Range("E17").Value = '  +2.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.814'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.48%  '
$ws.Range("D19").Value = '42.883.86'
$ws.Range("E19").Value = '  +4.47%  '
$ws.Range("E20").Value = '  +2.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.85%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.41'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.96%  '
$ws.Range("E24").Value = '  +12.82%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '230.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.31%  '
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.36'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.64%  '
$ws.Range("E29").Value = '  +1.73%  '
$ws.Range("E30").Value = '  +1.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '37.83'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +21.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '173.35'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.30'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.43%  '
$ws.Range("E34").Value = '  +2.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.38'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.88%  '
$ws.Range("E37").Value = '  +8.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.34'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.03%  '
$ws.Range("E39").Value = '  +14.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '13.07'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.36%  '
$ws.Range("E41").Value = '  +2.77%  '
$ws.Range("E42").Value = '  +2.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.200'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.63%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '60.04'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '105.44'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.61'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.04%  '
$ws.Range("E47").Value = '  +1.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.443'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +19.04%  '
$ws.Range("E49").Value = '  +1.28%  '
$ws.Range("E50").Value = '  +3.48%  '
$ws.Range("E51").Value = '  +0.37%  '
